$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; this shifts existing rows 16-67 down to 17-68
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with the new weekly data point
$ws.Cells.Item(16, 1).Value = 4
$ws.Cells.Item(16, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(16, 3).Value = "Los Lagos"
$ws.Cells.Item(16, 4).Value = 44998
$ws.Cells.Item(16, 5).Value = 10
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100101
$ws.Cells.Item(16, 8).Value = "Berries"
$ws.Cells.Item(16, 9).Value = 100101001
$ws.Cells.Item(16, 10).Value = "Arándano (blue)"
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 120
$ws.Cells.Item(16, 14).Value = 2300
$ws.Cells.Item(16, 15).Value = 2500
$ws.Cells.Item(16, 16).Value = 2400
$ws.Cells.Item(16, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(16, 19).Value = 1200
$ws.Cells.Item(16, 20).Value = 2
